$wb = $excel.ActiveWorkbook


# Sheet: Summary
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 0.01
$ws.Range("B6").Value = -55605.56227792593
$ws.Range("B7").Value = 10832798.46452648
$ws.Range("B8").Value = 21934850.53416022
$ws.Range("B10").Value = 4166178.291485272

# Sheet: Fed-in Capacity
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("K2").Value = 218.2111409464629
$ws.Range("L2").Value = 233.4357089658653
$ws.Range("M2").Value = 227.7528737187178
$ws.Range("N2").Value = 226.7777408199137
$ws.Range("O2").Value = 227.6097508526263
$ws.Range("P2").Value = 229.1091522321465
$ws.Range("K3").Value = 136.5653954824675
$ws.Range("L3").Value = 136.8385818314783
$ws.Range("M3").Value = 140.1317786452013
$ws.Range("N3").Value = 129.2864630672354
$ws.Range("O3").Value = 140.7160931202018
$ws.Range("P3").Value = 132.4654202195568
$ws.Range("Q3").Value = 138.9730561082392
$ws.Range("L4").Value = 133.8997196622861
$ws.Range("M4").Value = 137.8872852940284
$ws.Range("N4").Value = 126.6717389884799
$ws.Range("O4").Value = 137.5201250048898
$ws.Range("K5").Value = 208.3214547517356
$ws.Range("L5").Value = 221.1666787354324
$ws.Range("M5").Value = 214.1012135983212
$ws.Range("N5").Value = 212.9051825663978
$ws.Range("O5").Value = 214.5102874732696
$ws.Range("P5").Value = 217.9290633046618
$ws.Range("J6").Value = 122.1609165114202
$ws.Range("K6").Value = 129.8481964236088
$ws.Range("L6").Value = 127.8064787429921
$ws.Range("M6").Value = 129.5917403577182
$ws.Range("N6").Value = 118.4674613291454
$ws.Range("O6").Value = 130.8188201934042
$ws.Range("P6").Value = 124.5219861368073
$ws.Range("Q6").Value = 133.663080786811
$ws.Range("L7").Value = 128.7148261097669
$ws.Range("M7").Value = 132.4205420186164
$ws.Range("N7").Value = 121.3349826493852
$ws.Range("O7").Value = 132.590766746692
$ws.Range("P7").Value = 132.7088230120157
$ws.Range("O11").Value = 9.069265482343781
$ws.Range("Q11").Value = 80.64258426171736
$ws.Range("L13").Value = 47.39929435887774
$ws.Range("M13").Value = 46.68471506615327
$ws.Range("O41").Value = 9.069265482343809
$ws.Range("Q41").Value = 80.64258426171739
$ws.Range("L43").Value = 47.39929435887775
$ws.Range("M43").Value = 46.68471506615329

# Sheet: Unmet Demand
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("G2").Value = 415.2879682170718
$ws.Range("H2").Value = 339.3235460419765
$ws.Range("I2").Value = 209.9064962068211
$ws.Range("J2").Value = 10.69576364311372
$ws.Range("Q2").Value = 8.395781178313161
$ws.Range("R2").Value = 148.9413660216818
$ws.Range("S2").Value = 208.6835142066282
$ws.Range("T2").Value = 223.0311969618593
$ws.Range("U2").Value = 251.3444713639914
$ws.Range("G3").Value = 137.3356148920818
$ws.Range("H3").Value = 112.1591249337519
$ws.Range("I3").Value = 89.12455904281175
$ws.Range("J3").Value = 0
$ws.Range("R3").Value = 99.6672001609758
$ws.Range("S3").Value = 171.5363898835276
$ws.Range("T3").Value = 200.1328769967189
$ws.Range("U3").Value = 225.9408621947163
$ws.Range("G4").Value = 167.9843543597988
$ws.Range("H4").Value = 162.1682702466266
$ws.Range("I4").Value = 155.2512431493756
$ws.Range("J4").Value = 92.89079271141465
$ws.Range("K4").Value = 21.49978743611831
$ws.Range("P4").Value = 1.920177260816417
$ws.Range("Q4").Value = 85.60728995481479
$ws.Range("R4").Value = 176.9955073465141
$ws.Range("S4").Value = 223.9011423785076
$ws.Range("T4").Value = 227.9172826158254
$ws.Range("U4").Value = 286.3186679929276
$ws.Range("G5").Value = 415.2102214038314
$ws.Range("H5").Value = 338.5273214908784
$ws.Range("I5").Value = 206.909162189371
$ws.Range("J5").Value = 4.097100052852173
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 144.0576027644702
$ws.Range("S5").Value = 206.9118586999128
$ws.Range("T5").Value = 222.6908602868995
$ws.Range("U5").Value = 251.3382516189322
$ws.Range("G6").Value = 137.2940166796957
$ws.Range("H6").Value = 111.757373777286
$ws.Range("I6").Value = 87.69233988828964
$ws.Range("J6").Value = 0
$ws.Range("R6").Value = 97.08446202703757
$ws.Range("S6").Value = 170.7637213332854
$ws.Range("T6").Value = 199.9652070090748
$ws.Range("U6").Value = 225.9381254702172
$ws.Range("G7").Value = 167.9494798159088
$ws.Range("H7").Value = 161.8582038473137
$ws.Range("I7").Value = 154.2024705022108
$ws.Range("J7").Value = 90.42516245839145
$ws.Range("K7").Value = 17.44799951871617
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 82.6870224659892
$ws.Range("R7").Value = 175.4274210366964
$ws.Range("S7").Value = 223.2933741908973
$ws.Range("T7").Value = 227.7682732010226
$ws.Range("U7").Value = 286.3167657450791
$ws.Range("R11").Value = 67.46474657889402
$ws.Range("R12").Value = 56.57895837355434
$ws.Range("R41").Value = 67.46474657889404
$ws.Range("R42").Value = 56.57895837355435

# Sheet: Household Surplus
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B2").Value = 365891.4221755643
$ws.Range("B3").Value = 367442.3269559073

# Sheet: Costs and Revenues
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 64668.10444303645
$ws.Range("C2").Value = 65463.095369266
$ws.Range("D2").Value = 70000.02111109273
$ws.Range("K2").Value = 71864.55661388766
$ws.Range("M2").Value = 71864.55661388766
$ws.Range("P2").Value = 71864.55661388767
$ws.Range("B3").Value = 3540.182044145911
$ws.Range("C3").Value = 17610.60540304349
$ws.Range("D3").Value = 181359.9284513112
$ws.Range("B4").Value = 53328.91736755468
$ws.Range("C4").Value = 51037.36600445495
$ws.Range("B5").Value = 33707.69021107268
$ws.Range("C5").Value = 34129.29174257201
$ws.Range("B6").Value = -41045.72072275753
$ws.Range("C6").Value = -52384.9540799727
$ws.Range("D6").Value = -187269.9450821061
$ws.Range("E6").Value = -34791.82414645297
$ws.Range("F6").Value = 38750.23204350684
$ws.Range("G6").Value = 38750.23204350685
$ws.Range("H6").Value = 38750.23204350684
$ws.Range("I6").Value = 38750.23204350684
$ws.Range("J6").Value = 38750.23204350684
$ws.Range("K6").Value = 38750.23204350684
$ws.Range("L6").Value = 38750.23204350685
$ws.Range("M6").Value = 38750.23204350684
$ws.Range("N6").Value = 38750.23204350684
$ws.Range("O6").Value = 38750.23204350685
$ws.Range("P6").Value = 38750.23204350685

# Sheet: Installed Capacities
$ws = $wb.Worksheets.Item("Installed Capacities")
$ws.Range("B3").Value = 3.673862893242319
$ws.Range("C3").Value = 23.01338268678932

# Sheet: Added Capacities
$ws = $wb.Worksheets.Item("Added Capacities")
$ws.Range("B3").Value = 3.673862893242319
$ws.Range("C3").Value = 19.339519793547
$ws.Range("D3").Value = 211.7666554392049

# Sheet: PV Dispatch
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 0.01476929806328569
$ws.Range("H2").Value = 0.1512560737906247
$ws.Range("I2").Value = 0.5693933635848223
$ws.Range("J2").Value = 1.253525711498795
$ws.Range("K2").Value = 1.878710098517679
$ws.Range("L2").Value = 2.330706004121959
$ws.Range("M2").Value = 2.593359508554916
$ws.Range("N2").Value = 2.635322776677227
$ws.Range("O2").Value = 2.488460569060429
$ws.Range("P2").Value = 2.123843523123063
$ws.Range("Q2").Value = 1.594918036231644
$ws.Range("R2").Value = 0.9277519194678707
$ws.Range("S2").Value = 0.3365553796171231
$ws.Range("T2").Value = 0.06465260227203316
$ws.Range("U2").Value = 0.001181543845062855
$ws.Range("G3").Value = 0.007902271128860838
$ws.Range("H3").Value = 0.07631930274452442
$ws.Range("I3").Value = 0.2720738086033228
$ws.Range("J3").Value = 0.7465913262578567
$ws.Range("K3").Value = 1.276043491891533
$ws.Range("L3").Value = 1.715797948395859
$ws.Range("M3").Value = 2.002255276817064
$ws.Range("N3").Value = 2.05524901609789
$ws.Range("O3").Value = 1.880151324242605
$ws.Range("P3").Value = 1.508987194773435
$ws.Range("Q3").Value = 1.008717977782306
$ws.Range("R3").Value = 0.4906339916673424
$ws.Range("S3").Value = 0.1467812203102001
$ws.Range("T3").Value = 0.03185169810273293
$ws.Range("U3").Value = 0.000519886258477687
$ws.Range("G4").Value = 0.006624998659945164
$ws.Range("H4").Value = 0.05890226081296705
$ws.Range("I4").Value = 0.1992317778827147
$ws.Range("J4").Value = 0.4683874052581232
$ws.Range("K4").Value = 0.7697043897645381
$ws.Range("L4").Value = 0.9849566189522113
$ws.Range("M4").Value = 1.038498653576677
$ws.Range("N4").Value = 1.013805476753246
$ws.Range("O4").Value = 0.936413446952977
$ws.Range("P4").Value = 0.8012634742900949
$ws.Range("Q4").Value = 0.5547532968795902
$ws.Range("R4").Value = 0.2978840306553525
$ws.Range("S4").Value = 0.1154556584646807
$ws.Range("T4").Value = 0.02830681245612933
$ws.Range("U4").Value = 0.0003613635632697367
$ws.Range("G5").Value = 0.09251611130367558
$ws.Range("H5").Value = 0.9474806248887679
$ws.Range("I5").Value = 3.566727381034957
$ws.Range("J5").Value = 7.852189301760344
$ws.Range("K5").Value = 11.76839629324494
$ws.Range("L5").Value = 14.5997362345548
$ws.Range("M5").Value = 16.24501962895154
$ws.Range("N5").Value = 16.50788103019311
$ws.Range("O5").Value = 15.58792394841718
$ws.Range("P5").Value = 13.30393245060769
$ws.Range("Q5").Value = 9.990699214544804
$ws.Range("R5").Value = 5.811515176679516
$ws.Range("S5").Value = 2.108210886332509
$ws.Range("T5").Value = 0.4049892772318401
$ws.Range("U5").Value = 0.007401288904294046
$ws.Range("G6").Value = 0.0495004835149808
$ws.Range("H6").Value = 0.4780704592104726
$ws.Range("I6").Value = 1.704292963125436
$ws.Range("J6").Value = 4.676710155246499
$ws.Range("K6").Value = 7.993242550750212
$ws.Range("L6").Value = 10.74790103688213
$ws.Range("M6").Value = 12.54229356430018
$ws.Range("N6").Value = 12.87425075418792
$ws.Range("O6").Value = 11.77742425104019
$ws.Range("P6").Value = 9.452421277522957
$ws.Range("Q6").Value = 6.318693299210532
$ws.Range("R6").Value = 3.073372125605564
$ws.Range("S6").Value = 0.9194497705523841
$ws.Range("T6").Value = 0.1995216857467866
$ws.Range("U6").Value = 0.003256610757564528
$ws.Range("G7").Value = 0.04149954254994795
$ws.Range("H7").Value = 0.3689686601259011
$ws.Range("I7").Value = 1.248004425047526
$ws.Range("J7").Value = 2.93401765828132
$ws.Range("K7").Value = 4.821492307166679
$ws.Range("L7").Value = 6.169850171471354
$ws.Range("M7").Value = 6.505241928988658
$ws.Range("N7").Value = 6.350561815847949
$ws.Range("O7").Value = 5.865771705150827
$ws.Range("P7").Value = 5.019181037131885
$ws.Range("Q7").Value = 3.475020785705187
$ws.Range("R7").Value = 1.865970340473114
$ws.Range("S7").Value = 0.7232238460750018
$ws.Range("T7").Value = 0.1773162272588684
$ws.Range("U7").Value = 0.002263611411815345
$ws.Range("I11").Value = 50.5744059370523
$ws.Range("R11").Value = 82.40437136225566
$ws.Range("G12").Value = 0.70189203713141
$ws.Range("R12").Value = 43.57887577908879
$ws.Range("H13").Value = 5.231790602914802
$ws.Range("L13").Value = 87.48538192236055
$ws.Range("M13").Value = 92.24106888145178
$ws.Range("I41").Value = 50.57440593705229
$ws.Range("R41").Value = 82.40437136225565
$ws.Range("G42").Value = 0.7018920371314099
$ws.Range("R42").Value = 43.57887577908878
$ws.Range("H43").Value = 5.231790602914801
$ws.Range("L43").Value = 87.48538192236053
$ws.Range("M43").Value = 92.24106888145177
